$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables-specific")

# Select the entire row 111 (as Excel does before a row delete) and delete it,
# shifting all rows below it up by one.
$ws.Rows(111).Select()
$ws.Rows(111).Delete()

# Correct two typo'd values in the "flag_values" / "flag_meanings" cells for the
# qc_flag_quality_w* variables. The typo appears in three separate blocks
# (qc_flag_quality_wu/wv/wts), so fix it in all three places.
$fixedValues = "0, 1, 2, 3"
$fixedMeanings = "bad_data good_data good_for_reasearch suspect_data_good_for_general_use suspect_data_requires_further_checking_but_may_be_ok_for_general_use"

$ws.Range("C169").Value = $fixedMeanings
$ws.Range("C168").Value = $fixedValues
$ws.Range("C177").Value = $fixedMeanings
$ws.Range("C176").Value = $fixedValues
$ws.Range("C185").Value = $fixedMeanings
$ws.Range("C184").Value = $fixedValues

# Restore selection/scroll state similar to the post-edit session.
$ws.Range("A111").EntireRow.Select()
$ws.Application.ActiveWindow.ScrollRow = 81
